# Update the "Förändrad" (changed) date column C for rows 2-13
# from 2023-10-25 (serial 45224) to 2023-11-03 (serial 45233).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 13; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45224) {
        $cell.Value = 45233
    }
}
